$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 147
$ws.Range("F6").Value = 306
$ws.Range("F7").Value = 5446
$ws.Range("F9").Value = 7345
$ws.Range("F12").Value = 3753
$ws.Range("F13").Value = 61
$ws.Range("F20").Value = 94
$ws.Range("F22").Value = 3832
$ws.Range("F23").Value = 124
$ws.Range("F24").Value = 5121
$ws.Range("F25").Value = 433
$ws.Range("F26").Value = 2043
$ws.Range("F28").Value = 324
$ws.Range("F29").Value = 7572
$ws.Range("F30").Value = 28
$ws.Range("F32").Value = 2132
$ws.Range("F33").Value = 1319
$ws.Range("F34").Value = 147
$ws.Range("F35").Value = 1154
$ws.Range("F37").Value = 14
$ws.Range("F38").Value = 251
$ws.Range("F39").Value = 239
$ws.Range("F40").Value = 10
$ws.Range("F41").Value = 1172
$ws.Range("F42").Value = 1169
$ws.Range("F43").Value = 23
$ws.Range("F44").Value = 165
$ws.Range("F45").Value = 1295
$ws.Range("F46").Value = 1982
$ws.Range("F47").Value = 108
$ws.Range("F48").Value = 193
$ws.Range("F49").Value = 1203

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 39
$ws.Range("F3").Value = 14
$ws.Range("F4").Value = 144

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 538
$ws.Range("F3").Value = 714
$ws.Range("F4").Value = 61

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 147
$ws.Range("F5").Value = 538
$ws.Range("F6").Value = 714
$ws.Range("F7").Value = 61
$ws.Range("F8").Value = 306
$ws.Range("F9").Value = 5446
$ws.Range("F10").Value = 3753
$ws.Range("F11").Value = 61
$ws.Range("F17").Value = 94
$ws.Range("F18").Value = 14
$ws.Range("F19").Value = 144
$ws.Range("F21").Value = 3832
$ws.Range("F23").Value = 124
$ws.Range("F24").Value = 5121
$ws.Range("F25").Value = 433
$ws.Range("F26").Value = 2043
$ws.Range("F28").Value = 324
$ws.Range("F29").Value = 7572
$ws.Range("F30").Value = 28
$ws.Range("F32").Value = 2132
$ws.Range("F33").Value = 1319
$ws.Range("F34").Value = 147
$ws.Range("F35").Value = 1154
$ws.Range("F36").Value = 14
$ws.Range("F37").Value = 251
$ws.Range("F38").Value = 239
$ws.Range("F39").Value = 1172
$ws.Range("F40").Value = 1169
$ws.Range("F41").Value = 23
$ws.Range("F42").Value = 165
$ws.Range("F44").Value = 1295
$ws.Range("F46").Value = 1982
$ws.Range("F47").Value = 108
$ws.Range("F49").Value = 193
